# Adds part-number / stock info to the AR2ISS_MainCTL_HW BOM table.
# Mirrors commit "added all part no": a new "Spalte1" table column (G) plus
# a handful of previously-blank LCSC/Mouser part-number cells in columns E/F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AR2ISS_MainCTL_HW")
$lo = $ws.ListObjects.Item("AR2ISS_MainCTL_HW")

# --- 1. Grow the table by one column (name is set further down, in the
#        position that matches the original shared-string write order) ---
$newCol = $lo.ListColumns.Add()

# Give every data row in the new column the same "styled but empty" cell
# shape the rest of the table uses (s="1" General cells), by borrowing the
# format from an existing formatted column before filling in real values.
$fmtSource = $ws.Range("E2").Style
$ws.Range("G2:G75").Style = $fmtSource

# --- 2. Fill in newly-sourced part numbers / stock flags ------------------
# Order matches how the values were entered (bottom of sheet first), so the
# shared-string table comes out in the same append order as the real edit.

$ws.Range("F74").Value = "NTB0102DP-Q100H"          # U8  TXB0102DCT
$ws.Range("F72").Value = "Not Avaiable"             # U6  MCP2562-E-SN
$ws.Range("F64").Value = "700-MAX3232EUET"          # U1  MAX3232
$ws.Range("F65").Value = "TSR 1-2450E"              # U11 LME78_05-1.0
$ws.Range("F37").Value = "771-PMN48XP115"           # Q7  PMN48XP
$ws.Range("E14").Value = "C8598"                    # D26 D_Schottky
$ws.Range("E13").Value = "C96230"                   # D25 D31 D32 BZX84Cxx
$ws.Range("E15").Value = "C2128"                    # D3 D4 D7 D8 1N4148WS

$ws.Range("G11").Value = "xx"                       # D20 D_Zener
$ws.Range("G1").Value = "Spalte1"

$ws.Range("G19").Value = "stock"                    # J15 RJ45_LED

$ws.Range("F29").Value = "523-TSEH09SOL2RM8"        # J8  DB9_Female
$ws.Range("F34").Value = "81-BLM18KG102SN1D"        # L9 L26 L27 L
$ws.Range("E70").Value = "511-NUCLEO-H755ZI-Q"      # U2 NUCLEO144-H745ZI

# Remaining "xx" markers in the new column.
$ws.Range("G12").Value = "xx"                       # D21 D_TVS
$ws.Range("G18").Value = "xx"                       # J1 J2 Conn_01x06
$ws.Range("G20").Value = "xx"                       # J4 J17 Conn_01x02
$ws.Range("G21").Value = "xx"                       # J18 J19 Conn_02x17
$ws.Range("G22").Value = "xx"                       # J5 J22 J23 Conn_01x08
$ws.Range("G23").Value = "xx"                       # J24 Conn_02x10
$ws.Range("G24").Value = "xx"                       # J25 Conn_02x03
$ws.Range("G25").Value = "xx"                       # J27 Micro_SD_Card
$ws.Range("G27").Value = "xx"                       # J6 Conn_02x03
$ws.Range("G28").Value = "xx"                       # J7 Conn_01x03
$ws.Range("G30").Value = "xx"                       # J3 J9 ... Conn_01x02
$ws.Range("G32").Value = "xx"                       # K1 K2 K3 K4 FINDER-34.51

# --- 3. Selection / scroll position, to match the saved view -------------
$ws.Activate()
$ws.Range("E71").Select()
